# ----------------------------------------------------------------------------
# Rebuild the 'Code Smell' tracking table (header + 22 data rows, columns A-D).
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previous 18-row table (content + formatting) so it can be rebuilt.
$ws.Range("A1:D18").Clear()

# Full target table, written out row by row (A, B, C, D) for readability.
$table = @(
    @('Code Smell:', 'Location of Smell:', 'Refactored?', 'Reasoning'),
    @('Deficient Encapsulation', 'UserModel', 'Yes, partly', 'Some of the getters/setters needed to be set to private to ensure proper encapsulation.  Others could not be set to private as they are used in testing.'),
    @('Unutilized Abstraction', 'All Controllers/Impl classes', 'No', 'The classes in question are actually utilized in the front end and are there to serve a specific purpose'),
    @('Empty Class', 'EmptyPasswordException, InvalidCredentialsException, InvalidUserIdException', 'Yes', 'Added messages to all custom exceptions to fix empty class smells.'),
    @('Unnecessary Abstraction', 'EmptyPasswordException, InvalidCredentialsException, InvalidUserIdException', 'No', 'Custom exception classes, no need to refactor'),
    @('Unutilized Abstraction', 'WebSecurityConfig, TrelloCloneApplication, DatabaseConnectionTests, DashboardApplicationTests, ApplicationConstant', 'No', 'Classes in question are called in the front end/used for specific tasks'),
    @('Insufficient Modularization', 'TaskModel', 'No', 'Unnecessary to refactor, class contains mostly getters/setters and they are needed for functionality.'),
    @('Unnecessary Abstraction', 'DatabaseConnectionTests, DashboardApplicationTests', 'No', 'Unnecessary to refactor'),
    @('Deficient Encapsulation', 'ApplicationConstant', 'Yes', 'Change to enum and rename to Response to only contain message strings that will be used in the responses from the controller.'),
    @('Long Statement', 'TaskModel', 'Yes', 'Decompose the get methods t into various variables.'),
    @('Unutilized Abstraction', 'UnableTooAddBoardException', 'No', 'It’s an exception which is throws in the workspace service.'),
    @('Unutilized Abstraction', 'InvalidWorkspaceIdException', 'No', 'It’s an exception which is throws in the workspace service.'),
    @('Unutilized Abstraction', 'InvalidUserIdException', 'No', 'It’s an exception which is throws in the workspace and task service.'),
    @('Unutilized Abstraction', 'EmptyPasswordException', 'No', 'It’s an exception which is throws in the user service.'),
    @('Unutilized Abstraction', 'InvalidBoardIdException', 'No', 'It’s an exception which is throws in the workspace and board service.'),
    @('Unutilized Abstraction', 'InvalidTaskIdException', 'No', 'It’s an exception which is throws in the board and task service.'),
    @('Magic Number', 'BoardServiceImplTests', 'Yes', 'Global variables are introduced to replace the magic numbers'),
    @('Magic Number', 'TaskServiceImplTests', 'Yes', 'Global variables are introduced to replace the magic numbers'),
    @('Magic Number', 'UserServiceImplTests', 'Yes', 'Global variables are introduced to replace the magic numbers'),
    @('Long Statement', 'TaskController', 'Yes', 'Introduce variables to decompose the long statements.'),
    @('Unnecessary Abstraction', 'TaskStatusEnum', 'No', 'Unnecessary to refactor'),
    @('Unnecessary Abstraction', 'ApplicationConstant', 'No', 'Unnecessary to refactor'),
)

# Excel's shared-string table records each unique string the first time it is
# written, so cells are poked in the same first-seen order as the reference
# workbook rather than strictly row-by-row: first one cell per distinct string,
# then any repeats (which just reuse the already-registered string).
$seen = @{}
$firstPass = @()
$secondPass = @()
for ($r = 0; $r -lt $table.Count; $r++) {
    for ($c = 0; $c -lt 4; $c++) {
        $value = $table[$r][$c]
        $entry = @{ Row = ($r + 1); Col = ($c + 1); Value = $value }
        if ($seen.ContainsKey($value)) {
            $secondPass += $entry
        } else {
            $seen[$value] = $true
            $firstPass += $entry
        }
    }
}

foreach ($entry in $firstPass) {
    $ws.Cells.Item($entry.Row, $entry.Col).Value = $entry.Value
}
foreach ($entry in $secondPass) {
    $ws.Cells.Item($entry.Row, $entry.Col).Value = $entry.Value
}

# Rows 11-22 (the black-font rows in the reference workbook) get an explicit
# black font color, matching the pre-existing 's="1"' cell style.
$ws.Range("A11:D22").Font.Color = 0

# Match the saved selection/active cell from the reference workbook.
$ws.Range("A4").Select()
